$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.004.36"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.260.42"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.66"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.45"
$ws.Range("E6").Value = "  -2.97%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.258.64"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.510"
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("E12").Value = "  -1.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000242"
$ws.Range("E13").Value = "  -2.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.87"
$ws.Range("E14").Value = "  -1.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.798.49"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.121"
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.258.77"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.085.67"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.71"
$ws.Range("E19").Value = "  -1.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "470.56"
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.77"
$ws.Range("E21").Value = "  -3.52%  "
$ws.Range("E22").Value = "  -1.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.82"
$ws.Range("E23").Value = "  -2.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.49"
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.81"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -0.67%  "
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.05"
$ws.Range("E29").Value = "  -2.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.91"
$ws.Range("E30").Value = "  -2.34%  "
$ws.Range("E31").Value = "  -1.89%  "
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("E33").Value = "  -3.03%  "
$ws.Range("E34").Value = "  -3.67%  "
$ws.Range("E35").Value = "  -1.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.88"
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.77"
$ws.Range("E37").Value = "  -1.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0717"
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.070.02"
$ws.Range("E40").Value = "  +2.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "421.26"
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.117"
$ws.Range("E42").Value = "  +6.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.18"
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.64"
$ws.Range("E44").Value = "  -4.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.257"
$ws.Range("E45").Value = "  -3.17%  "
$ws.Range("E46").Value = "  -1.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.69"
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("E51").Value = "  -2.11%  "

# Row 48/49: Monero and Arweave swap positions with updated values
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.72"
$ws.Range("E48").Value = "  +6.59%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "126.81"
$ws.Range("E49").Value = "  +3.31%  "
